$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 52633940
$ws.Range("I106").Value = 90910840
$ws.Range("J106").Value = 3212
$ws.Range("K106").Value = 90910840
$ws.Range("L106").Value = 3212
$ws.Range("M106").Value = -90910209
$ws.Range("N106").Value = -4474
$ws.Range("H132").Value = 5850436
$ws.Range("I132").Value = 2620.7778
$ws.Range("K132").Value = 7862.3334
$ws.Range("M132").Value = -5332.3334
$ws.Range("H137").Value = 21443812
$ws.Range("I137").Value = 5682825.5
$ws.Range("J137").Value = 48116252
$ws.Range("K137").Value = 17048476.5
$ws.Range("L137").Value = 144348756
$ws.Range("M137").Value = -17045926.5
$ws.Range("N137").Value = -144353856
$ws.Range("H141").Value = 1468.2858
$ws.Range("I141").Value = 1468.2858
$ws.Range("K141").Value = 4404.857400000001
$ws.Range("M141").Value = 775.1425999999992

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 36390
$ws.Range("I2").Value = 59600
$ws.Range("J2").Value = 1575
$ws.Range("K2").Value = 59600
$ws.Range("L2").Value = 1575
$ws.Range("M2").Value = -59487
$ws.Range("N2").Value = -1801
$ws.Range("H32").Value = 3702.95
$ws.Range("I32").Value = 2567.6707
$ws.Range("J32").Value = 10136.2
$ws.Range("K32").Value = 2567.6707
$ws.Range("L32").Value = 10136.2
$ws.Range("M32").Value = -2280.6707
$ws.Range("N32").Value = -10710.2
$ws.Range("H45").Value = 334270.56
$ws.Range("I45").Value = 715107.4
$ws.Range("K45").Value = 715107.4
$ws.Range("M45").Value = -714730.4
$ws.Range("H61").Value = 2792339.5
$ws.Range("I61").Value = 1544169.1
$ws.Range("J61").Value = 6536850.5
$ws.Range("K61").Value = 1544169.1
$ws.Range("L61").Value = 6536850.5
$ws.Range("M61").Value = -1543957.1
$ws.Range("N61").Value = -6537274.5
$ws.Range("H116").Value = 36390
$ws.Range("I116").Value = 59600
$ws.Range("J116").Value = 1575
$ws.Range("K116").Value = 59600
$ws.Range("L116").Value = 1575
$ws.Range("M116").Value = -57306
$ws.Range("N116").Value = -6163
$ws.Range("H136").Value = 2792339.5
$ws.Range("I136").Value = 1544169.1
$ws.Range("J136").Value = 6536850.5
$ws.Range("K136").Value = 4632507.300000001
$ws.Range("L136").Value = 19610551.5
$ws.Range("M136").Value = -4629957.300000001
$ws.Range("N136").Value = -19615651.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 36390
$ws.Range("I3").Value = 59600
$ws.Range("J3").Value = 1575
$ws.Range("K3").Value = 59600
$ws.Range("L3").Value = 1575
$ws.Range("M3").Value = -59486
$ws.Range("N3").Value = -1803
$ws.Range("H80").Value = 276.6316
$ws.Range("I80").Value = 104.333336
$ws.Range("J80").Value = 356.15384
$ws.Range("K80").Value = 104.333336
$ws.Range("L80").Value = 356.15384
$ws.Range("M80").Value = 893.666664
$ws.Range("N80").Value = -2352.15384
$ws.Range("H83").Value = 276.6316
$ws.Range("I83").Value = 104.333336
$ws.Range("J83").Value = 356.15384
$ws.Range("K83").Value = 521.66668
$ws.Range("L83").Value = 1780.7692
$ws.Range("M83").Value = 4470.33332
$ws.Range("N83").Value = -11764.7692
$ws.Range("H86").Value = 1989
$ws.Range("I86").Value = 1998.9796
$ws.Range("K86").Value = 1998.9796
$ws.Range("M86").Value = -875.9795999999999
$ws.Range("H89").Value = 1989
$ws.Range("I89").Value = 1998.9796
$ws.Range("K89").Value = 9994.897999999999
$ws.Range("M89").Value = -4378.897999999999
$ws.Range("H105").Value = 1924.6666
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 1864.4
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 1864.4
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5358.4
$ws.Range("H107").Value = 1012.875
$ws.Range("I107").Value = 1020.6
$ws.Range("K107").Value = 1020.6
$ws.Range("M107").Value = 899.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1049.75
$ws.Range("I16").Value = 1049.75
$ws.Range("K16").Value = 1049.75
$ws.Range("M16").Value = -762.75
$ws.Range("H31").Value = 2481726.2
$ws.Range("I31").Value = 1345597
$ws.Range("K31").Value = 1345597
$ws.Range("M31").Value = -1345302
$ws.Range("H34").Value = 2481726.2
$ws.Range("I34").Value = 1345597
$ws.Range("K34").Value = 1345597
$ws.Range("M34").Value = -1345395
$ws.Range("H113").Value = 1049.75
$ws.Range("I113").Value = 1049.75
$ws.Range("K113").Value = 1049.75
$ws.Range("M113").Value = 1120.25
$ws.Range("H132").Value = 2364.8276
$ws.Range("I132").Value = 1662.2273
$ws.Range("J132").Value = 4573
$ws.Range("K132").Value = 4986.6819
$ws.Range("L132").Value = 13719
$ws.Range("M132").Value = -2456.6819
$ws.Range("N132").Value = -18779
$ws.Range("H134").Value = 1487495
$ws.Range("I134").Value = 6006.5713
$ws.Range("J134").Value = 6672704.5
$ws.Range("K134").Value = 18019.7139
$ws.Range("L134").Value = 20018113.5
$ws.Range("M134").Value = -15484.7139
$ws.Range("N134").Value = -20023183.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3130.074
$ws.Range("I22").Value = 2018.4615
$ws.Range("J22").Value = 4162.2856
$ws.Range("K22").Value = 2018.4615
$ws.Range("L22").Value = 4162.2856
$ws.Range("M22").Value = -1723.4615
$ws.Range("N22").Value = -4752.2856
$ws.Range("H27").Value = 3130.074
$ws.Range("I27").Value = 2018.4615
$ws.Range("J27").Value = 4162.2856
$ws.Range("K27").Value = 2018.4615
$ws.Range("L27").Value = 4162.2856
$ws.Range("M27").Value = -1911.4615
$ws.Range("N27").Value = -4376.2856
$ws.Range("H122").Value = 11648700
$ws.Range("I122").Value = 1332368.8
$ws.Range("K122").Value = 3997106.4
$ws.Range("M122").Value = -3994656.4
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10914.565
$ws.Range("I136").Value = 7205.278
$ws.Range("J136").Value = 24268
$ws.Range("K136").Value = 21615.834
$ws.Range("L136").Value = 72804
$ws.Range("M136").Value = -19065.834
$ws.Range("N136").Value = -77904

# --- Special case: row 128 on LTW loses its N cell entirely (was a computed profit field that no longer applies) ---
$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("N128").ClearContents()
